$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Nome", "Sobrenome", "Idade")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $ws.Cells.Item(1, $c + 1)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.Font.Italic = $true
}

$data = @(
    @("Daniel", "Galleazzo", 19),
    @("Paulo", "Galleazzo", 21),
    @("Júlia", "Zanon", 20),
    @("Sandra", "Galleazzo", 51),
    @("Antônio", "Galleazzo", 81),
    @("Maria", "Galleazzo", 80),
    @("Raissa", "AnticristoSDD", 666)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

$ws.Range("A2:C2").Font.Name = "Calibri"

$ws.Columns.Item(1).ColumnWidth = 12
$ws.Columns.Item(2).ColumnWidth = 13.28515625
